$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Insert the new "2000-09 spinup" worksheet right before "CW3M c118 2010"
# ---------------------------------------------------------------------
$before = $wb.Worksheets.Item("CW3M c118 2010")
$spin = $wb.Worksheets.Add($before)
$spin.Name = "2000-09 spinup"

$src = $wb.Worksheets.Item("2010-18")

# Column widths / header formatting to match sheet "2010-18"
$spin.Columns("B").ColumnWidth = 28.88671875

# Copy the header row (row 1) values from "2010-18"
$src.Range("A1:S1").Copy()
$spin.Range("A1").PasteSpecial(-4163)

# Re-apply the header cell formatting
$spin.Rows(1).RowHeight = 129.6
$spin.Range("A1:S1").WrapText = $true
$spin.Range("C1").HorizontalAlignment = -4108
$spin.Range("D1:N1").NumberFormat = "0.00"
$spin.Range("O1:P1").NumberFormat = "0"
$spin.Range("R1").NumberFormat = "0.000000"

# Data rows 2 & 3
$spin.Range("A2").Value = "CW3M"
$spin.Range("B2").Value = "Baseline 2000-09 5/9/21 spinup"
$spin.Range("C2").Value = "2000-09"
$spin.Range("D2").Value = 572.4274934
$spin.Range("E2").Value = 1951.2097047
$spin.Range("F2").Value = 6.0977103000000001
$spin.Range("G2").Value = 195.51971589999999
$spin.Range("H2").Value = 0
$spin.Range("I2").Value = 7.3807704000000003
$spin.Range("J2").Value = 0
$spin.Range("K2").Value = 545.82672409999998
$spin.Range("L2").Value = 92.905869899999999
$spin.Range("M2").Value = 1490.1102661999998
$spin.Range("N2").Value = 603.64529430000005
$spin.Range("O2").Value = 16485.282812500001
$spin.Range("P2").Value = 1985.1201415999999
$spin.Range("Q2").Value = -0.14723959999999991
$spin.Range("R2").Value = -0.00032939999999999998

$spin.Range("A3").Value = "CW3M"
$spin.Range("B3").Value = "Baseline 2000-09 C393 spinup"
$spin.Range("C3").Value = "2000-09"
$spin.Range("D3").Value = 571.75505380000004
$spin.Range("E3").Value = 1951.2097047
$spin.Range("F3").Value = 5.8274805000000001
$spin.Range("G3").Value = 195.51971589999999
$spin.Range("H3").Value = 0
$spin.Range("I3").Value = 7.3943439999999994
$spin.Range("J3").Value = 0
$spin.Range("K3").Value = 545.72043469999994
$spin.Range("L3").Value = 93.221763699999997
$spin.Range("M3").Value = 1489.4764649000001
$spin.Range("N3").Value = 603.14073799999994
$spin.Range("O3").Value = 15579.5007324
$spin.Range("P3").Value = 1985.1201415999999
$spin.Range("Q3").Value = -0.14689749999999976
$spin.Range("R3").Value = -0.00033019999999999989

# Number formats for the two data rows
$spin.Range("D2:N3").NumberFormat = "0.00"
$spin.Range("O2:P3").NumberFormat = "0"
$spin.Range("O3").NumberFormat = "0"
$spin.Range("Q2:Q3").NumberFormat = "0.00"
$spin.Range("R2:R3").NumberFormat = "0.000000"

# ---------------------------------------------------------------------
# 2. Add row 12 to "2010-18" (copy of row 11, new scenario label)
# ---------------------------------------------------------------------
$src.Range("A11:R11").Copy()
$src.Range("A12").PasteSpecial(-4163)
$src.Range("B12").Value = "Baseline 2010-18 C393"

$src.Range("D12:N12").NumberFormat = "0.00"
$src.Range("O12:P12").NumberFormat = "0"
$src.Range("Q12").NumberFormat = "0.00"
$src.Range("R12").NumberFormat = "0.000000"

# ---------------------------------------------------------------------
# 3. Selection / active sheet bookkeeping
# ---------------------------------------------------------------------
$src.Range("A1:XFD1").Select()
$spin.Activate()
$spin.Range("O3").Select()
